# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.083.13'
$ws.Range('E2').Value = '  -2.10%  '

# Row 3
$ws.Range('D3').Value = '3.480.44'
$ws.Range('E3').Value = '  -0.88%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.65'
$ws.Range('E5').Value = '  -2.57%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.42'
$ws.Range('E6').Value = '  -4.12%  '

# Row 7
$ws.Range('D7').Value = '3.479.63'
$ws.Range('E7').Value = '  -0.93%  '

# Row 8
$ws.Range('E8').Value = '  +0.07%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.491'
$ws.Range('E9').Value = '  -2.77%  '

# Row 10
$ws.Range('E10').Value = '  -5.03%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.384'
$ws.Range('E12').Value = '  -4.49%  '

# Row 13
$ws.Range('D13').Value = '4.070.58'
$ws.Range('E13').Value = '  -0.70%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000181'
$ws.Range('E14').Value = '  -6.32%  '

# Row 15
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.61'
$ws.Range('E15').Value = '  -6.57%  '

# Row 16
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '3.451.12'
$ws.Range('E16').Value = '  -1.95%  '

# Row 17
$ws.Range('E17').Value = '  -1.16%  '

# Row 18
$ws.Range('D18').Value = '65.061.38'
$ws.Range('E18').Value = '  -1.94%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.74'
$ws.Range('E19').Value = '  -8.16%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '5.77'
$ws.Range('E20').Value = '  -4.95%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.96'
$ws.Range('E21').Value = '  -3.99%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '389.46'
$ws.Range('E22').Value = '  -7.32%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.557'
$ws.Range('E23').Value = '  -4.65%  '

# Row 24
$ws.Range('B24').Value = 'Dai'
$ws.Range('C24').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.03%  '

# Row 25
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.54'
$ws.Range('E25').Value = '  -5.40%  '

# Row 26
$ws.Range('B26').Value = 'WrappedeETH'
$ws.Range('C26').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D26').Value = '3.618.77'
$ws.Range('E26').Value = '  -0.99%  '

# Row 27
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.78'
$ws.Range('E27').Value = '  +0.79%  '

# Row 28
$ws.Range('E28').Value = '  -1.68%  '

# Row 29
$ws.Range('E29').Value = '  -0.19%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.41'
$ws.Range('E30').Value = '  -4.01%  '

# Row 31
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.22'
$ws.Range('E31').Value = '  -9.17%  '

# Row 32
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.13'
$ws.Range('E32').Value = '  -8.49%  '

# Row 33
$ws.Range('D33').Value = '3.497.59'
$ws.Range('E33').Value = '  -0.56%  '

# Row 34
$ws.Range('E34').Value = '  -0.05%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.144'
$ws.Range('E35').Value = '  -6.02%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.06'
$ws.Range('E36').Value = '  -4.27%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '169.79'
$ws.Range('E37').Value = '  -2.48%  '

# Row 38
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.20'
$ws.Range('E38').Value = '  -9.36%  '

# Row 39
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.83'
$ws.Range('E39').Value = '  -8.30%  '

# Row 40
$ws.Range('E40').Value = '  -9.22%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.73'
$ws.Range('E41').Value = '  -8.52%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0778'
$ws.Range('E42').Value = '  -3.02%  '

# Row 43
$ws.Range('E43').Value = '  -4.57%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.67'
$ws.Range('E44').Value = '  -6.16%  '

# Row 45
$ws.Range('E45').Value = '  +0.06%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '25.15'
$ws.Range('E46').Value = '  +10.58%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.35'
$ws.Range('E47').Value = '  -11.66%  '

# Row 48
$ws.Range('E48').Value = '  +5.10%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '6.70'
$ws.Range('E50').Value = '  -4.26%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.07'
$ws.Range('E51').Value = '  -10.78%  '
